$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.354.05'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.687.86'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('D4').Value = '''1.011'
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').Value = '''218.71'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = '''0.5449'
$ws.Range('E6').Value = '  +3.93%  '
$ws.Range('D7').Value = '''1.010'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('D8').Value = '''0.2744'
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').Value = '''0.06449'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '''21.96'
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('D11').Value = '''0.07691'
$ws.Range('D12').Value = '1.685.93'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '''4.526'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').Value = '''0.000008384'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').Value = '''65.15'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '26.396.85'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '''4.936'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '''10.96'
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').Value = '''191.40'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').Value = '''6.251'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '''1.011'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').Value = '''149.48'
$ws.Range('E24').Value = '  +2.92%  '
$ws.Range('D25').Value = '''0.1322'
$ws.Range('E25').Value = '  +6.64%  '
$ws.Range('D26').Value = '''7.876'
$ws.Range('E26').Value = '  +2.56%  '
$ws.Range('D27').Value = '''15.74'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').Value = '''0.06359'
$ws.Range('E28').Value = '  -4.79%  '
$ws.Range('D29').Value = '''1.402'
$ws.Range('E29').Value = '  +4.12%  '
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').Value = '''3.595'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = '''3.577'
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').Value = '''1.682'
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('D35').Value = '''0.6149'
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').Value = '1.117.66'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('D40').Value = '''0.01631'
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('D41').Value = '''0.8769'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D43').Value = '''101.79'
$ws.Range('E43').Value = '  +0.89%  '
$ws.Range('D44').Value = '1.837.61'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '''57.44'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = '''0.00000000107'
$ws.Range('E46').Value = '  -7.59%  '
$ws.Range('D47').Value = '''1.016'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('D48').Value = '''8.188'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').Value = '''0.4304'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '''6.087'
$ws.Range('E51').Value = '  +1.15%  '
